$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row text updates (I1, L1) ---
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# --- Move the "Status as of July 4, 2025" column from AA to AF, ---
# --- and insert 5 new "No. of Sites ..." header columns in AA:AE ---
$ws.Range("AA1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("AF1").Value = "Status as of July 4, 2025"

# New headers AA1:AE1 - copy the bold/bordered style from an existing
# styled header cell (H1 uses style index 1) then set their text.
$ws.Range("H1").Copy()
$ws.Range("AA1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# --- Remove the now-unused "-" placeholder values in columns I and L ---
# for data rows 2-8 (cells become fully empty, not just value-cleared)
$ws.Range("I2:I8").ClearContents()
$ws.Range("L2:L8").ClearContents()

# --- Update the dropdown data validation to target the relocated column ---
$ws.Range("AF2:AF8").Validation.Delete()
$ws.Range("AF2:AF8").Validation.Add(3, 1, 1, "=DropdownOptions!$A$1:$A$7")
$ws.Range("AF2:AF8").Validation.IgnoreBlank = $true
$ws.Range("AF2:AF8").Validation.InCellDropdown = $false

$ws.Range("AA2:AA8").Validation.Delete()
